# Auto-generated script to apply scraped diff changes to Rafflesia_Profits workbook
# Columns H-N are plain cached numeric values (no formulas) per-row leve market data.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 5 (Leve Item ID 5503)
$ws.Range("H5").Value = 151
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
# Row 20 (Leve Item ID 1965)
$ws.Range("H20").Value = 13498
$ws.Range("I20").Value = 13498
$ws.Range("K20").Value = 13498
$ws.Range("M20").Value = -13268
# Row 35 (Leve Item ID 1965)
$ws.Range("H35").Value = 13498
$ws.Range("I35").Value = 13498
$ws.Range("K35").Value = 13498
$ws.Range("M35").Value = -13119
# Row 38 (Leve Item ID 4599)
$ws.Range("H38").Value = 11271.917
$ws.Range("J38").Value = 19800
$ws.Range("L38").Value = 59400
$ws.Range("N38").Value = -60144
# Row 43 (Leve Item ID 5472)
$ws.Range("H43").Value = 1900
$ws.Range("I43").Value = 1000
$ws.Range("J43").Value = 3250
$ws.Range("K43").Value = 1000
$ws.Range("L43").Value = 3250
$ws.Range("M43").Value = -931
$ws.Range("N43").Value = -3388
# Row 58 (Leve Item ID 4606)
$ws.Range("H58").Value = 2256.875
$ws.Range("I58").Value = 343.33334
$ws.Range("J58").Value = 3405
$ws.Range("K58").Value = 1030.00002
$ws.Range("L58").Value = 10215
$ws.Range("M58").Value = -880.0000199999999
$ws.Range("N58").Value = -10515
# Row 61 (Leve Item ID 4604)
$ws.Range("H61").Value = 483.33334
$ws.Range("I61").Value = 225
$ws.Range("K61").Value = 675
$ws.Range("M61").Value = -503
# Row 101 (Leve Item ID 19884)
$ws.Range("H101").Value = 529.3333
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 5249.75
$ws.Range("I32").Value = 4666.3335
$ws.Range("K32").Value = 4666.3335
$ws.Range("M32").Value = -4379.3335
# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 1314.3334
$ws.Range("I45").Value = 1314.3334
$ws.Range("K45").Value = 1314.3334
$ws.Range("M45").Value = -937.3334
# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 2575.7144
$ws.Range("I122").Value = 1979.3529
$ws.Range("K122").Value = 5938.0587
$ws.Range("M122").Value = -3488.0587

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 22 (Leve Item ID 5092)
$ws.Range("H22").Value = 346.5
$ws.Range("I22").Value = 340.33334
$ws.Range("K22").Value = 340.33334
$ws.Range("M22").Value = -167.33334
# Row 50 (Leve Item ID 27159)
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
# Row 107 (Leve Item ID 27706)
$ws.Range("H107").Value = 4583
$ws.Range("I107").Value = 9999
$ws.Range("K107").Value = 9999
$ws.Range("M107").Value = -8079

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 10 (Leve Item ID 1997)
$ws.Range("H10").Value = 5165
$ws.Range("I10").Value = 247.5
$ws.Range("J10").Value = 15000
$ws.Range("K10").Value = 247.5
$ws.Range("L10").Value = 15000
$ws.Range("M10").Value = -108.5
$ws.Range("N10").Value = -15278
# Row 16 (Leve Item ID 27691)
$ws.Range("H16").Value = 3404.1428
$ws.Range("I16").Value = 3603.6667
$ws.Range("J16").Value = 3254.5
$ws.Range("K16").Value = 3603.6667
$ws.Range("L16").Value = 3254.5
$ws.Range("M16").Value = -3316.6667
$ws.Range("N16").Value = -3828.5
# Row 22 (Leve Item ID 5367)
$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 500
$ws.Range("K22").Value = 500
$ws.Range("M22").Value = -150
# Row 23 (Leve Item ID 2703)
$ws.Range("H23").Value = 400
$ws.Range("J23").Value = 400
$ws.Range("L23").Value = 400
$ws.Range("N23").Value = -880
# Row 27 (Leve Item ID 2703)
$ws.Range("H27").Value = 400
$ws.Range("J27").Value = 400
$ws.Range("L27").Value = 400
$ws.Range("N27").Value = -784
# Row 45 (Leve Item ID 2026)
$ws.Range("H45").Value = 4000
$ws.Range("I45").Value = 4000
$ws.Range("K45").Value = 4000
$ws.Range("M45").Value = -3407
# Row 68 (Leve Item ID 10611)
$ws.Range("H68").Value = 50000
$ws.Range("J68").Value = 50000
$ws.Range("L68").Value = 50000
$ws.Range("N68").Value = -51498
# Row 71 (Leve Item ID 10611)
$ws.Range("H71").Value = 50000
$ws.Range("J71").Value = 50000
$ws.Range("L71").Value = 150000
$ws.Range("N71").Value = -157488
# Row 99 (Leve Item ID 36198)
$ws.Range("H99").Value = 5218.6665
$ws.Range("I99").Value = 4937.3335
$ws.Range("K99").Value = 4937.3335
$ws.Range("M99").Value = -3439.3335
# Row 105 (Leve Item ID 19928)
$ws.Range("H105").Value = 1007.1429
$ws.Range("I105").Value = 990
$ws.Range("K105").Value = 990
$ws.Range("M105").Value = 757
# Row 113 (Leve Item ID 27691)
$ws.Range("H113").Value = 3404.1428
$ws.Range("I113").Value = 3603.6667
$ws.Range("J113").Value = 3254.5
$ws.Range("K113").Value = 3603.6667
$ws.Range("L113").Value = 3254.5
$ws.Range("M113").Value = -1433.6667
$ws.Range("N113").Value = -7594.5
# Row 126 (Leve Item ID 36198)
$ws.Range("H126").Value = 5218.6665
$ws.Range("I126").Value = 4937.3335
$ws.Range("K126").Value = 14812.0005
$ws.Range("M126").Value = -12342.0005
# Row 127 (Leve Item ID 35351)
$ws.Range("H127").Value = 27777.777
$ws.Range("I127").Value = 20000
$ws.Range("J127").Value = 30000
$ws.Range("K127").Value = 20000
$ws.Range("L127").Value = 30000
$ws.Range("M127").Value = -15040
$ws.Range("N127").Value = -39920

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 17 (Leve Item ID 4640)
$ws.Range("H17").Value = 52.5
$ws.Range("I17").Value = 50.666668
$ws.Range("J17").Value = 55.25
$ws.Range("K17").Value = 152.000004
$ws.Range("L17").Value = 165.75
$ws.Range("M17").Value = 16.99999600000001
$ws.Range("N17").Value = -503.75
# Row 34 (Leve Item ID 4749)
$ws.Range("H34").Value = 1371.3125
$ws.Range("I34").Value = 513.8
$ws.Range("J34").Value = 1761.091
$ws.Range("K34").Value = 1541.4
$ws.Range("L34").Value = 5283.272999999999
$ws.Range("M34").Value = -1457.4
$ws.Range("N34").Value = -5451.272999999999
# Row 38 (Leve Item ID 4860)
$ws.Range("H38").Value = 2439.5715
$ws.Range("J38").Value = 4549.5
$ws.Range("L38").Value = 13648.5
$ws.Range("N38").Value = -14342.5
# Row 39 (Leve Item ID 4712)
$ws.Range("H39").Value = 1225
$ws.Range("I39").Value = 450
$ws.Range("K39").Value = 1350
$ws.Range("M39").Value = -1056
# Row 55 (Leve Item ID 4733)
$ws.Range("H55").Value = 3139.2144
$ws.Range("I55").Value = 1620
$ws.Range("J55").Value = 3983.2222
$ws.Range("K55").Value = 4860
$ws.Range("L55").Value = 11949.6666
$ws.Range("M55").Value = -4683
$ws.Range("N55").Value = -12303.6666

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 9 (Leve Item ID 1683)
$ws.Range("H9").Value = 13098
$ws.Range("I9").Value = 2745
$ws.Range("J9").Value = 20000
$ws.Range("K9").Value = 2745
$ws.Range("L9").Value = 20000
$ws.Range("M9").Value = -2575
$ws.Range("N9").Value = -20340
# Row 40 (Leve Item ID 4113)
$ws.Range("H40").Value = 1269
$ws.Range("J40").Value = 1269
$ws.Range("L40").Value = 1269
$ws.Range("N40").Value = -1571
# Row 43 (Leve Item ID 4218)
$ws.Range("H43").Value = 10466.5
$ws.Range("I43").Value = 5699.75
$ws.Range("J43").Value = 20000
$ws.Range("K43").Value = 5699.75
$ws.Range("L43").Value = 20000
$ws.Range("M43").Value = -5548.75
$ws.Range("N43").Value = -20302
# Row 44 (Leve Item ID 4143)
$ws.Range("H44").Value = 13000.25
$ws.Range("I44").Value = 8000
$ws.Range("J44").Value = 14667
$ws.Range("K44").Value = 8000
$ws.Range("L44").Value = 14667
$ws.Range("M44").Value = -7404
$ws.Range("N44").Value = -15859
# Row 48 (Leve Item ID 4337)
$ws.Range("H48").Value = 19625.25
$ws.Range("I48").Value = 7500
$ws.Range("J48").Value = 23667
$ws.Range("K48").Value = 7500
$ws.Range("L48").Value = 23667
$ws.Range("M48").Value = -7015
$ws.Range("N48").Value = -24637

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 20 (Leve Item ID 4308)
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
# Row 30 (Leve Item ID 1688)
$ws.Range("H30").Value = 6521
$ws.Range("I30").Value = 878.5714
$ws.Range("J30").Value = 46018
$ws.Range("K30").Value = 878.5714
$ws.Range("L30").Value = 46018
$ws.Range("M30").Value = -770.5714
$ws.Range("N30").Value = -46234
# Row 35 (Leve Item ID 1697)
$ws.Range("H35").Value = 1735.8334
$ws.Range("I35").Value = 1276
$ws.Range("J35").Value = 4035
$ws.Range("K35").Value = 1276
$ws.Range("L35").Value = 4035
$ws.Range("M35").Value = -940
$ws.Range("N35").Value = -4707
# Row 61 (Leve Item ID 27740)
$ws.Range("H61").Value = 850
$ws.Range("I61").Value = 850
$ws.Range("K61").Value = 850
$ws.Range("M61").Value = -648
# Row 113 (Leve Item ID 27740)
$ws.Range("H113").Value = 850
$ws.Range("I113").Value = 850
$ws.Range("K113").Value = 850
$ws.Range("M113").Value = 1320
# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 11518
$ws.Range("I136").Value = 3124
$ws.Range("K136").Value = 9372
$ws.Range("M136").Value = -6822

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81 (Leve Item ID 12596)
$ws.Range("H81").Value = 7283.3335
$ws.Range("I81").Value = 10425
$ws.Range("K81").Value = 20850
$ws.Range("M81").Value = -19789
# Row 84 (Leve Item ID 12596)
$ws.Range("H84").Value = 7283.3335
$ws.Range("I84").Value = 10425
$ws.Range("K84").Value = 104250
$ws.Range("M84").Value = -98946
# Row 100 (Leve Item ID 19981)
$ws.Range("H100").Value = 1011
$ws.Range("I100").Value = 1096.1666
$ws.Range("K100").Value = 2192.3332
$ws.Range("M100").Value = -1651.3332

